# "User data with filename fix"
# - Notes sheet: update "Units of measure" line to reflect the real unit.
# - Data sheet: fill in the two data rows that were previously blank.

$wb = $excel.ActiveWorkbook

$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A3").Value = 'Units of measure: constant 2015 US$'

$data = $wb.Worksheets.Item("Data")

$data.Range("A2").Value = 'bilateral-unspecified'
$data.Range("B2").Value = 'Bilateral, unspecified'
$data.Range("C2").Value = 2015
$data.Range("D2").Value = 217650000

$data.Range("A3").Value = 'north-central-america'
$data.Range("B3").Value = 'North & Central America, regional'
$data.Range("C3").Value = 2015
$data.Range("D3").Value = 7390000
